$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New column C: "update" date header + date value -----------------------
$ws.Range("C1").Value = "עדכון"
$ws.Range("C2").NumberFormat = "d-mmm"
$ws.Range("C2").Value = Get-Date -Year 2026 -Month 1 -Day 14 -Hour 0 -Minute 0 -Second 0

# --- Append the new ranking rows (rows 14-66) -------------------------------
$newRows = @(
    @("יובל סטרוזר", 1),
    @("הילס שולויס", 1),
    @("אורי שטרנברג", 1),
    @("איתי הראל", 1),
    @("אורי שטרנברג", 1),
    @("תומר ששון", 1),
    @("אן מרש", 1),
    @("איתי הראל", 1),
    @("יובל סטרוזר", 1),
    @("ליהי בראל", 1),
    @("יהלי דוייב", 1),
    @("ירון גלפנד", 1),
    @("גלי זליג", 1),
    @("ליהי בראל", 1),
    @("איתי בסטקר", 1),
    @("דפנה ברגשטיין", 1),
    @("ליאם דיין ", 1),
    @("יהלי דוייב", 1),
    @("תומר ששון", 1),
    @("יולי יערי תליו", 1),
    @("אן מרש", 1),
    @("יהלי גודר", 1),
    @("ירון גלפנד", 1),
    @("איתי הראל", 1),
    @("יולי יערי תליו", 1),
    @("דפנה ברגשטיין", 1),
    @("נועם מילר", 1),
    @("ליהי בראל", 1),
    @("הילה שולויס", 1),
    @("גלי זליג", 1),
    @("אן מרש", 1),
    @("דפנה ברגשטיין", 1),
    @("מעיין סטרוזר", 1),
    @("נועם מילר", 1),
    @("יולי יערי תליו", 1),
    @("ליאם דיין ", 1),
    @("תאיו ורד", 1),
    @("איתי הראל", 6),
    @("ליאם מלכה", 6),
    @("תומר ששון", 6),
    @("איתי הראל", 6),
    @("איתי הראל", 6),
    @("יהלי דוייב", 6),
    @("יולי יערי תליו", 6),
    @("ליאם מלכה", 6),
    @("תומר ששון", 6),
    @("הילס שולויס", 6),
    @("אן מרש", 6),
    @("ליהי בראל", 6),
    @("תאיו ורד", 6),
    @("דפנה ברגשטיין", 6),
    @("ירון גלפנד", 6),
    @("מעיין סטרוזר", 6)
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

$ws.Range("C3").Select()
